$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -5.8998
$ws.Range("E6").Value = 12.46180000000001
$ws.Range("E7").Value = 12.39169999999999
$ws.Range("D8").Value = -8.918799999999985
$ws.Range("E8").Value = 12.63129999999999
$ws.Range("B12").Value = 5.1275
$ws.Range("D12").Value = -8.247800000000003
$ws.Range("D14").Value = -8.790900000000001
$ws.Range("E19").Value = 12.79809999999999
$ws.Range("E21").Value = 12.91349999999999
$ws.Range("D22").Value = -8.182599999999994
$ws.Range("E24").Value = 12.93979999999999
